# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (these columns store prices as plain text).
$textCells = @('D5','D6','D8','D10','D15','D16','D19','D20','D23','D24','D25','D29','D35','D39','D40','D41','D43','D44','D46','D48','D50')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row-by-row cell updates
# Row 2
$ws.Range('D2').Value = '26.125.59'
$ws.Range('E2').Value = '  -2.37%  '

# Row 3
$ws.Range('D3').Value = '1.572.66'
$ws.Range('E3').Value = '  -1.93%  '

# Row 4
$ws.Range('E4').Value = '  -0.44%  '

# Row 5
$ws.Range('D5').Value = '208.23'
$ws.Range('E5').Value = '  -1.78%  '

# Row 6
$ws.Range('D6').Value = '0.499'
$ws.Range('E6').Value = '  -2.96%  '

# Row 7
$ws.Range('E7').Value = '  -0.41%  '

# Row 8
$ws.Range('D8').Value = '0.0609'
$ws.Range('E8').Value = '  -1.71%  '

# Row 9
$ws.Range('E9').Value = '  -1.31%  '

# Row 10
$ws.Range('D10').Value = '19.56'
$ws.Range('E10').Value = '  -0.82%  '

# Row 11
$ws.Range('E11').Value = '  -0.40%  '

# Row 12
$ws.Range('D12').Value = '1.792.28'

# Row 13
$ws.Range('D13').Value = '1.576.94'
$ws.Range('E13').Value = '  -2.09%  '

# Row 14
$ws.Range('E14').Value = '  -0.71%  '

# Row 15
$ws.Range('D15').Value = '0.514'
$ws.Range('E15').Value = '  -2.37%  '

# Row 16
$ws.Range('D16').Value = '64.28'
$ws.Range('E16').Value = '  -1.24%  '

# Row 17
$ws.Range('D17').Value = '26.113.49'
$ws.Range('E17').Value = '  -2.28%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0725'
$ws.Range('E18').Value = '  -2.25%  '

# Row 19
$ws.Range('D19').Value = '7.30'
$ws.Range('E19').Value = '  +2.15%  '

# Row 20
$ws.Range('D20').Value = '207.79'
$ws.Range('E20').Value = '  -1.29%  '

# Row 21
$ws.Range('E21').Value = '  -0.43%  '

# Row 22
$ws.Range('E22').Value = '  -1.62%  '

# Row 23
$ws.Range('D23').Value = '2.18'
$ws.Range('E23').Value = '  -2.83%  '

# Row 24
$ws.Range('D24').Value = '8.81'
$ws.Range('E24').Value = '  -3.00%  '

# Row 25
$ws.Range('D25').Value = '143.45'
$ws.Range('E25').Value = '  -0.25%  '

# Row 26
$ws.Range('E26').Value = '  -0.38%  '

# Row 27
$ws.Range('E27').Value = '  -1.88%  '

# Row 28
$ws.Range('E28').Value = '  -2.09%  '

# Row 29
$ws.Range('D29').Value = '15.21'
$ws.Range('E29').Value = '  -1.20%  '

# Row 30
$ws.Range('E30').Value = '  -0.46%  '

# Row 31
$ws.Range('E31').Value = '  -1.61%  '

# Row 32
$ws.Range('E32').Value = '  -2.14%  '

# Row 33
$ws.Range('E33').Value = '  +0.28%  '

# Row 34
$ws.Range('D34').Value = '1.276.65'
$ws.Range('E34').Value = '  -1.43%  '

# Row 35
$ws.Range('D35').Value = '0.610'
$ws.Range('E35').Value = '  +2.97%  '

# Row 36
$ws.Range('E36').Value = '  -1.57%  '

# Row 37
$ws.Range('E37').Value = '  -1.29%  '

# Row 38
$ws.Range('E38').Value = '  -2.92%  '

# Row 39
$ws.Range('D39').Value = '1.09'
$ws.Range('E39').Value = '  -10.06%  '

# Row 40
$ws.Range('D40').Value = '0.810'
$ws.Range('E40').Value = '  -2.68%  '

# Row 41
$ws.Range('D41').Value = '5.55'
$ws.Range('E41').Value = '  +1.88%  '

# Row 42
$ws.Range('E42').Value = '  -2.84%  '

# Row 43
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '62.38'
$ws.Range('E43').Value = '  -1.24%  '

# Row 44
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '0.761'
$ws.Range('E44').Value = '  -2.62%  '

# Row 45
$ws.Range('D45').Value = '1.706.41'
$ws.Range('E45').Value = '  -1.93%  '

# Row 46
$ws.Range('D46').Value = '88.92'
$ws.Range('E46').Value = '  -1.73%  '

# Row 47
$ws.Range('E47').Value = '  -3.25%  '

# Row 48
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.0999'
$ws.Range('E48').Value = '  -2.40%  '

# Row 49
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₇0998'
$ws.Range('E49').Value = '  -3.32%  '

# Row 50
$ws.Range('D50').Value = '0.0506'
$ws.Range('E50').Value = '  -1.38%  '

# Row 51
$ws.Range('E51').Value = '  -0.31%  '
